# Apply the 06:44:15 scrape update to the "horarios-141" workbook.
# Sheet 1 = LP1912, Sheet 2 = LP1912-215, Sheet 3 = 6203-6173

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("A2").Value = "Última actualización: 06:44:15"
$ws1.Range("A3").Value = "Total filas: 53"

# Make room for the 16 newly-scraped rows that get merged (by arrival time)
# into the existing list, pushing rows 29..42 down.
$ws1.Range("A29:A44").EntireRow.Insert()

# Final, fully-merged contents for rows 29..58.
$sheet1Rows = @(
    @(29, "06:44:15", "07:00", "14_ABASTO", 16, "LP1912"),
    @(30, "06:44:15", "07:01", "16_SANTA ANA", 17, "LP1912"),
    @(31, "05:31:23", "07:04", "23_HERNANDEZ", 93, "LP1912"),
    @(32, "05:31:23", "07:05", "15_ABASTO", 94, "LP1912"),
    @(33, "06:44:15", "07:05", "23_HERNANDEZ", 21, "LP1912"),
    @(34, "05:31:23", "07:07", "225_GOMEZ", 96, "LP1912"),
    @(35, "05:31:23", "07:11", "215A_EL PATO", 100, "LP1912"),
    @(36, "05:31:23", "07:15", "11_ETCHEVERRY", 104, "LP1912"),
    @(37, "06:44:15", "07:16", "11_ETCHEVERRY", 32, "LP1912"),
    @(38, "06:44:15", "07:17", "16_SANTA ANA", 33, "LP1912"),
    @(39, "05:31:23", "07:21", "26_HERNANDEZ", 110, "LP1912"),
    @(40, "05:31:23", "07:23", "10_OLMOS", 112, "LP1912"),
    @(41, "06:44:15", "07:25", "10_OLMOS", 41, "LP1912"),
    @(42, "06:01:37", "07:31", "11_ETCHEVERRY", 90, "LP1912"),
    @(43, "06:01:37", "07:31", "16_SANTA ANA", 90, "LP1912"),
    @(44, "06:01:37", "07:32", "84_COLONIA URQUIZA-ESC 49", 91, "LP1912"),
    @(45, "06:44:15", "07:32", "11_ETCHEVERRY", 48, "LP1912"),
    @(46, "06:01:37", "07:36", "27_EL RETIRO", 95, "LP1912"),
    @(47, "06:44:15", "07:37", "27_EL RETIRO", 53, "LP1912"),
    @(48, "06:01:37", "07:39", "10_OLMOS", 98, "LP1912"),
    @(49, "06:01:37", "07:47", "14_ABASTO", 106, "LP1912"),
    @(50, "06:44:15", "07:48", "14_ABASTO", 64, "LP1912"),
    @(51, "06:01:37", "07:51", "215D_EL PATO", 110, "LP1912"),
    @(52, "06:44:15", "08:04", "23_HERNANDEZ", 80, "LP1912"),
    @(53, "06:44:15", "08:12", "15_ABASTO", 88, "LP1912"),
    @(54, "06:44:15", "08:21", "26_HERNANDEZ", 97, "LP1912"),
    @(55, "06:44:15", "08:23", "16_P MOR-SANTA ANA", 99, "LP1912"),
    @(56, "06:44:15", "08:23", "215B_EL PATO", 99, "LP1912"),
    @(57, "06:44:15", "08:27", "84_COLONIA URQUIZA-ESC 49", 103, "LP1912"),
    @(58, "06:44:15", "08:42", "81_EL PELIGRO", 118, "LP1912")
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("A2").Value = "Última actualización: 06:44:15"
$ws2.Range("A3").Value = "Total filas: 9"

$ws2.Cells.Item(14, 1).Value = "06:44:15"
$ws2.Cells.Item(14, 2).Value = "08:23"
$ws2.Cells.Item(14, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(14, 4).Value = 99
$ws2.Cells.Item(14, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("A2").Value = "Última actualización: 06:44:15"
$ws3.Range("A3").Value = "Total filas: 11"

$sheet3Rows = @(
    @(14, "06:44:15", "07:36", "215A_LA PLATA", 52, "L6173"),
    @(15, "06:44:15", "08:07", "215C_LA PLATA", 83, "L6203"),
    @(16, "06:44:15", "08:31", "215A_LA PLATA", 107, "L6173")
)

foreach ($row in $sheet3Rows) {
    $r = $row[0]
    $ws3.Cells.Item($r, 1).Value = $row[1]
    $ws3.Cells.Item($r, 2).Value = $row[2]
    $ws3.Cells.Item($r, 3).Value = $row[3]
    $ws3.Cells.Item($r, 4).Value = $row[4]
    $ws3.Cells.Item($r, 5).Value = $row[5]
}
